# Adding result classes to the image importer
# - Insert a new column before the old "Fixes all?" column (J) to hold a new
#   "Operation" classification column. Excel will shift the old J/K columns to
#   K/L and carry formatting along automatically.
# - Populate the new "Operation" column with the per-image results.
# - Re-select the cell that was active when the workbook was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J; existing J:K data shifts right to K:L.
$ws.Columns("J:J").Insert()

# Header for the new column.
$ws.Range("J2").Value = "Operation"

# Per-row "Operation" values for the rows that had a classification.
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("J9").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("J14").Value = 2
$ws.Range("J16").Value = "?"
$ws.Range("J17").Value = "?"
$ws.Range("J18").Value = "?"

# Match the new column's width with its neighbours (H:I).
$ws.Columns("J:J").ColumnWidth = $ws.Columns("I:I").ColumnWidth

# Restore the active-cell selection saved with the workbook.
[void]$ws.Range("H4").Select()
